$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the two new worksheets at the end of the workbook
# ---------------------------------------------------------------------------
$srcMass2 = $wb.Worksheets.Item("E_Glass_linear_full_2mass")

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsFbTemp = $wb.Worksheets.Add($null, $lastSheet)
$wsFbTemp.Name = "E_Glass_linear_full_mass_fbtemp"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsHybrid = $wb.Worksheets.Add($null, $lastSheet)
$wsHybrid.Name = "multiReactionCharringHybrid"


# ---------------------------------------------------------------------------
# 2. Populate "E_Glass_linear_full_mass_fbtemp" - comparison table like the
#    other "2mass"/"full" sheets (copy formatting from the sibling sheet,
#    then overwrite with this sheet's own data).
# ---------------------------------------------------------------------------
$srcMass2.Range("A1:D19").Copy()
$wsFbTemp.Range("A1:D19").PasteSpecial(-4122) # xlPasteFormats
$wsFbTemp.Application.CutCopyMode = $false

$wsFbTemp.Rows.Item(1).RowHeight = 30.75
for ($r = 2; $r -le 18; $r++) {
    $wsFbTemp.Rows.Item($r).RowHeight = 15.75
}
$wsFbTemp.Range("A1:D19").ColumnWidth = 20.7109375

# Header row
$wsFbTemp.Range("A1").Value = "Input parameters"
$wsFbTemp.Range("B1").Value = "Virtual inputs"
$wsFbTemp.Range("C1").Value = "Hybrid method"
$wsFbTemp.Range("D1").Value = "Hybrid error (%)"

# Row 2 - A (1/s)
$wsFbTemp.Range("A2").Value = "A (1/s)"
$wsFbTemp.Range("B2").Value = 2340000000000
$wsFbTemp.Range("C2").Value = 3430000000000
$wsFbTemp.Range("D2").Formula = "=ABS(C2-B2)/B2"

# Row 3 - E (J/mol)
$wsFbTemp.Range("A3").Value = "E (J/mol)"
$wsFbTemp.Range("B3").Value = 181000
$wsFbTemp.Range("C3").Value = 183000
$wsFbTemp.Range("D3").Formula = "=ABS(C3-B3)/B3"

# Row 4 - n
$wsFbTemp.Range("A4").Value = "n"
$wsFbTemp.Range("B4").Value = 1
$wsFbTemp.Range("C4").Value = 1.013
$wsFbTemp.Range("D4").Formula = "=ABS(C4-B4)/B4"

# Row 5 - Qpyro (J/kg)
$wsFbTemp.Range("A5").Value = "Qpyro (J/kg)"
$wsFbTemp.Range("B5").Value = 100000
$wsFbTemp.Range("C5").Value = 99847.5
$wsFbTemp.Range("D5").Formula = "=ABS(C5-B5)/B5"

# Row 6 - kv_a (W/m/K)
$wsFbTemp.Range("A6").Value = "kv_a (W/m/K)"
$wsFbTemp.Range("B6").Value = 0.312
$wsFbTemp.Range("C6").Value = 0.313
$wsFbTemp.Range("D6").Formula = "=ABS(C6-B6)/B6"

# Row 7 - kv_b (W/m/K^2)
$wsFbTemp.Range("A7").Value = "kv_b (W/m/K^2)"
$wsFbTemp.Range("B7").Value = 0.00004405
$wsFbTemp.Range("C7").Value = 0.0000424
$wsFbTemp.Range("D7").Formula = "=ABS(C7-B7)/B7"

# Row 8 - Cpv_a (J/kg/K)
$wsFbTemp.Range("A8").Value = "Cpv_a (J/kg/K)"
$wsFbTemp.Range("B8").Value = 1080
$wsFbTemp.Range("C8").Value = 1082.2
$wsFbTemp.Range("D8").Formula = "=ABS(C8-B8)/B8"

# Row 9 - Cpv_b (J/kg/K^2)
$wsFbTemp.Range("A9").Value = "Cpv_b (J/kg/K^2)"
$wsFbTemp.Range("B9").Value = 0.0452
$wsFbTemp.Range("C9").Value = 0.04576
$wsFbTemp.Range("D9").Formula = "=ABS(C9-B9)/B9"

# Row 10 - kc_a (W/m/K)
$wsFbTemp.Range("A10").Value = "kc_a (W/m/K)"
$wsFbTemp.Range("B10").Value = 0.0949
$wsFbTemp.Range("C10").Value = 0.0834
$wsFbTemp.Range("D10").Formula = "=ABS(C10-B10)/B10"

# Row 11 - kc_b (W/m/K^2)
$wsFbTemp.Range("A11").Value = "kc_b (W/m/K^2)"
$wsFbTemp.Range("B11").Value = 0.000283
$wsFbTemp.Range("C11").Value = 0.00029999999999999997
$wsFbTemp.Range("D11").Formula = "=ABS(C11-B11)/B11"

# Row 12 - Cpc_a (J/kg/K)
$wsFbTemp.Range("A12").Value = "Cpc_a (J/kg/K)"
$wsFbTemp.Range("B12").Value = 1041
$wsFbTemp.Range("C12").Value = 1045
$wsFbTemp.Range("D12").Formula = "=ABS(C12-B12)/B12"

# Row 13 - Cpc_b (J/kg/K^2)
$wsFbTemp.Range("A13").Value = "Cpc_b (J/kg/K^2)"
$wsFbTemp.Range("B13").Value = 0.259
$wsFbTemp.Range("C13").Value = 0.2565
$wsFbTemp.Range("D13").Formula = "=ABS(C13-B13)/B13"

# Row 14 - rho v (kg/m^3) - no hybrid value
$wsFbTemp.Range("A14").Value = "ρv (kg/m^3)"
$wsFbTemp.Range("B14").Value = 1683

# Row 15 - rho c (kg/m^3) - no hybrid value
$wsFbTemp.Range("A15").Value = "ρc (kg/m^3)"
$wsFbTemp.Range("B15").Value = 1235

# Row 16 - epsilon v
$wsFbTemp.Range("A16").Value = "εv"
$wsFbTemp.Range("B16").Value = 0.94
$wsFbTemp.Range("C16").Value = 0.941716
$wsFbTemp.Range("D16").Formula = "=ABS(C16-B16)/B16"

# Row 17 - epsilon c - no hybrid value
$wsFbTemp.Range("A17").Value = "εc"
$wsFbTemp.Range("B17").Value = 0.94

# Row 18 - Kv (1/m)
$wsFbTemp.Range("A18").Value = "Kv (1/m)"
$wsFbTemp.Range("B18").Value = 10000
$wsFbTemp.Range("C18").Value = 10246.6
$wsFbTemp.Range("D18").Formula = "=ABS(C18-B18)/B18"

# Row 19 - footer label
$wsFbTemp.Range("A19").Value = "ρv, ρc, εc fixed"

$wsFbTemp.Range("D29").Select()

Write-Host "Sheets after add:"
foreach ($s in $wb.Worksheets) {
    Write-Host "  $($s.Name)"
}
